# feat: add 2022-Q4 data
#
# - Insert a new worksheet named "2022-Q4" right after "总计", pushing the
#   existing "2022-Q3" sheet to 3rd position.
# - On "总计", the former 2022-Q3 totals row moves down to row 3 and row 2 is
#   overwritten with the new 2022-Q4 totals (10 holdings, 0.66 billion yuan).
# - The new "2022-Q4" sheet gets the per-fund holdings table (mirrors the
#   layout already used on the "2022-Q3" sheet).

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet between "总计" and "2022-Q3".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total, 1)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) "总计" sheet: push the old row (2022-Q3 totals) down to row 3, then
#    overwrite row 2 with the new 2022-Q4 totals.
# ---------------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 5
$total.Range("D3").Value = 1.05

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.66

# ---------------------------------------------------------------------------
# 3) Populate the new "2022-Q4" sheet with the fund holdings table, reusing
#    the header / row-number style ("s=2", same as used on "总计") from the
#    existing sheets.
# ---------------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4data = @(
    @(0, "166011", "中欧盛世成长混合（LOF）A", "4.97", "92.56", "4.39", "0.2182", 3),
    @(1, "001306", "中欧永裕混合A",           "3.71", "92.63", "4.38", "0.1625", 3),
    @(2, "008445", "融通产业趋势先锋股票",       "3.01", "91.03", "4.65", "0.1400", 2),
    @(3, "580006", "东吴新经济混合A",           "1.16", "89.15", "5.25", "0.0609", 4),
    @(4, "004233", "中欧盛世成长混合（LOF）C", "0.52", "92.56", "4.39", "0.0228", 3),
    @(5, "012617", "东吴新经济混合C",           "0.42", "89.15", "5.25", "0.0220", 4),
    @(6, "001307", "中欧永裕混合C",             "0.29", "92.63", "4.38", "0.0127", 3),
    @(7, "002409", "华夏新活力灵活配置混合A",    "0.12", "77.58", "9.40", "0.0113", 2),
    @(8, "001888", "中欧盛世成长混合（LOF）E", "0.25", "92.56", "4.39", "0.0110", 3),
    @(9, "002410", "华夏新活力灵活配置混合C",    "0.00", "77.58", "9.40", $null,   2)
)

$r = 2
foreach ($row in $q4data) {
    $q4.Range("A$r").Value = $row[0]

    # Fund code (B) and the numeric-looking measures (D/E/F/G) are stored as
    # TEXT in the source data (mirrors how the "2022-Q3" sheet is laid out) -
    # a leading "'" keeps Excel from coercing them to numbers, and resetting
    # the Style back to Normal afterwards drops the quote-prefix cell style
    # so the cell stays plain (no explicit "s" attribute), same as the diff.
    $q4.Range("B$r").Formula = "'" + $row[1]
    $q4.Range("B$r").Style = "Normal"

    $q4.Range("C$r").Value = $row[2]

    $q4.Range("D$r").Formula = "'" + $row[3]
    $q4.Range("D$r").Style = "Normal"
    $q4.Range("E$r").Formula = "'" + $row[4]
    $q4.Range("E$r").Style = "Normal"
    $q4.Range("F$r").Formula = "'" + $row[5]
    $q4.Range("F$r").Style = "Normal"
    if ($row[6] -eq $null) {
        $q4.Range("G$r").Value = 0
    } else {
        $q4.Range("G$r").Formula = "'" + $row[6]
        $q4.Range("G$r").Style = "Normal"
    }
    $q4.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Row-number column (A) and header row (B1:H1) reuse the bold/boxed style ("s=2")
# that the "总计" sheet's A-column already carries.
$total.Range("A2").Copy()
$q4.Range("A2:A11").PasteSpecial(-4122)
$q4.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Restore "总计" as the active sheet/tab (adding a sheet made the new one
#    active).
# ---------------------------------------------------------------------------
$total.Activate()
$total.Range("A1").Select()
